$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates per diff: Price (D) and Volume(1h) (E) columns, rows 2-51.
# D-column values that parse as plain numbers need an explicit Text number
# format first, otherwise Excel auto-converts them to numeric cells and
# the literal text (e.g. trailing zeros like "0.04850") would be lost.

$ws.Range('D2').Value = '26.010.80'
$ws.Range('E2').Value = '  +0.59%  '
$ws.Range('D3').Value = '1.641.46'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.73'
$ws.Range('E5').Value = '  +0.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5092'
$ws.Range('E6').Value = '  +1.45%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2562'
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06354'
$ws.Range('E9').Value = '  -0.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.53'
$ws.Range('E10').Value = '  +0.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07760'
$ws.Range('E11').Value = '  -0.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.281'
$ws.Range('E12').Value = '  +0.48%  '
$ws.Range('D13').Value = '1.641.53'
$ws.Range('E13').Value = '  -0.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5433'
$ws.Range('E14').Value = '  +0.42%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '64.21'
$ws.Range('E15').Value = '  -0.56%  '
$ws.Range('D16').Value = '0.0₅7712'
$ws.Range('E16').Value = '  -1.80%  '
$ws.Range('D17').Value = '26.034.89'
$ws.Range('E17').Value = '  +0.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.001'
$ws.Range('E18').Value = '  -0.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '197.39'
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.411'
$ws.Range('E20').Value = '  +0.82%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.911'
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.025'
$ws.Range('E22').Value = '  +1.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.003'
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.856'
$ws.Range('E24').Value = '  -0.79%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '140.84'
$ws.Range('E25').Value = '  +0.70%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1193'
$ws.Range('E26').Value = '  +4.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.812'
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.58'
$ws.Range('E28').Value = '  -0.56%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.233'
$ws.Range('E29').Value = '  -0.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.04850'
$ws.Range('E30').Value = '  -0.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.247'
$ws.Range('E31').Value = '  -0.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.158'
$ws.Range('E32').Value = '  -0.84%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.522'
$ws.Range('E33').Value = '  -0.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.361'
$ws.Range('E34').Value = '  -0.21%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.8968'
$ws.Range('E35').Value = '  +1.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.575'
$ws.Range('E36').Value = '  -0.91%  '
$ws.Range('D37').Value = '1.139.26'
$ws.Range('E37').Value = '  +0.66%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5449'
$ws.Range('E38').Value = '  -1.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01559'
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.001'
$ws.Range('E40').Value = '  -0.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.541'
$ws.Range('E41').Value = '  -0.56%  '
$ws.Range('E42').Value = '  +6.58%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8094'
$ws.Range('E43').Value = '  -0.89%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.36'
$ws.Range('E44').Value = '  +0.09%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.384'
$ws.Range('E45').Value = '  -5.14%  '
$ws.Range('D46').Value = '1.780.20'
$ws.Range('E46').Value = '  +0.28%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4520'
$ws.Range('E47').Value = '  +0.22%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9984'
$ws.Range('E48').Value = '  -0.88%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.75'
$ws.Range('E49').Value = '  -0.54%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05052'
$ws.Range('E50').Value = '  -0.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.004'
$ws.Range('E51').Value = '  -0.18%  '
